$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (before Doveton), shifting existing rows down
$ws.Rows("3:3").Insert()

$ws.Range("A3").Value = "Caulfield"
$ws.Range("B3").Value = "Metro Train - Frankston line"
$ws.Range("C3").Value = "30/12/20 4:30pm-17:00pm"
$ws.Range("D3").Value = "Case caught train from Caulfield to Cheltenham"

# Insert a new row at position 11 (before what is now Wyanga Winery row), shifting existing rows down
$ws.Rows("11:11").Insert()

$ws.Range("A11").Value = "Lakes Entrance"
$ws.Range("B11").Value = "V/Line bus - Lakes Entrance to Bairnsdale"
$ws.Range("C11").Value = "30/12/2020 11:55am-12:30pm"
$ws.Range("D11").Value = "Case caught the 11:55am bus from Lakes Entrance"
